$wb = $excel.ActiveWorkbook

# --- "Overview" sheet: the per-locale status columns (E = zh-cn, F = de-de)
# show "Ready for handoff" for every row; the report regeneration moved the
# files into translation, so the status text becomes "In Translation".
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("E2:F4").Value = "In Translation"

# --- Per-locale detail sheets: column C ("Status") holds the same status
# string for each row.
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("C2:C4").Value = "In Translation"

$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("C2:C4").Value = "In Translation"

# The shorter replacement text means the status columns can be narrower.
# ColumnWidth is specified in characters (using the workbook's default
# font); 12.5 is the value that maps to the report's target stored column
# width for these cells.
$ws1.Columns.Item(5).ColumnWidth = 12.5
$ws1.Columns.Item(6).ColumnWidth = 12.5
$ws2.Columns.Item(3).ColumnWidth = 12.5
$ws3.Columns.Item(3).ColumnWidth = 12.5

Write-Output "Report regenerated: status -> 'In Translation'; status columns resized."
